$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 107, shifting existing rows 107-124 down to 108-125
$ws.Rows("107:107").Insert()

# Populate the newly inserted row 107 with the new weekly data point
$ws.Cells.Item(107, 1).Value = 9
$ws.Cells.Item(107, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(107, 3).Value = "Metropolitana"
$ws.Cells.Item(107, 4).Value = [DateTime]"2023-03-20"
$ws.Cells.Item(107, 5).Value = 13
$ws.Cells.Item(107, 6).Value = 100114007
$ws.Cells.Item(107, 7).Value = "Jengibre"
$ws.Cells.Item(107, 8).Value = "Sin especificar"
$ws.Cells.Item(107, 9).Value = "Primera"
$ws.Cells.Item(107, 10).Value = 610
$ws.Cells.Item(107, 11).Value = 17000
$ws.Cells.Item(107, 12).Value = 18000
$ws.Cells.Item(107, 13).Value = 17500
$ws.Cells.Item(107, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(107, 15).Value = "Perú"
$ws.Cells.Item(107, 16).Value = 1346
$ws.Cells.Item(107, 17).Value = 13
$ws.Cells.Item(107, 18).Value = "Hortaliza"
